# Update "想去人数" (interested-people count) figures in column F,
# reflecting refreshed data generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 4265
$ws1.Range("F7").Value  = 3393
$ws1.Range("F11").Value = 272
$ws1.Range("F13").Value = 1243
$ws1.Range("F17").Value = 238
$ws1.Range("F19").Value = 9463
$ws1.Range("F23").Value = 797
$ws1.Range("F25").Value = 822
$ws1.Range("F31").Value = 221
$ws1.Range("F33").Value = 4766
$ws1.Range("F35").Value = 987
$ws1.Range("F36").Value = 120

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 8611
$ws3.Range("F4").Value = 1480

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 8611
$ws4.Range("F5").Value  = 1480
$ws4.Range("F7").Value  = 4265
$ws4.Range("F10").Value = 3393
$ws4.Range("F14").Value = 272
$ws4.Range("F20").Value = 1243
$ws4.Range("F25").Value = 238
$ws4.Range("F27").Value = 9463
$ws4.Range("F32").Value = 797
$ws4.Range("F34").Value = 822
$ws4.Range("F39").Value = 221
$ws4.Range("F42").Value = 4766
$ws4.Range("F44").Value = 987
